$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.402881622314453
$ws.Range("B1").Value = 1.554775953292847
$ws.Range("C1").Value = 1.625626564025879
$ws.Range("D1").Value = 1.493659853935242
$ws.Range("E1").Value = 1.245018482208252
